$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.485.30"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3
$ws.Range("D3").Value = "1.645.73"
$ws.Range("E3").Value = "  -1.25%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'212.71"
$ws.Range("E5").Value = "  -1.34%  "

# Row 6
$ws.Range("E6").Value = "  +3.68%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").Value = "'23.50"
$ws.Range("E8").Value = "  -2.86%  "

# Row 9
$ws.Range("D9").Value = "'0.258"
$ws.Range("E9").Value = "  -2.32%  "

# Row 10
$ws.Range("D10").Value = "'0.0612"
$ws.Range("E10").Value = "  -1.44%  "

# Row 11
$ws.Range("D11").Value = "'0.0892"

# Row 12
$ws.Range("D12").Value = "1.877.40"
$ws.Range("E12").Value = "  -1.32%  "

# Row 13
$ws.Range("D13").Value = "1.649.85"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("D14").Value = "'0.590"
$ws.Range("E14").Value = "  +4.07%  "

# Row 15
$ws.Range("D15").Value = "'4.05"
$ws.Range("E15").Value = "  -2.11%  "

# Row 16
$ws.Range("D16").Value = "'64.54"
$ws.Range("E16").Value = "  -2.98%  "

# Row 17
$ws.Range("D17").Value = "27.446.95"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18
$ws.Range("D18").Value = "'231.15"
$ws.Range("E18").Value = "  -4.59%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.98%  "

# Row 20
$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = "  -1.73%  "

# Row 21
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").Value = "'4.35"
$ws.Range("E22").Value = "  -3.59%  "

# Row 23
$ws.Range("D23").Value = "'9.73"
$ws.Range("E23").Value = "  +3.78%  "

# Row 24
$ws.Range("E24").Value = "  -1.41%  "

# Row 25
$ws.Range("D25").Value = "'147.76"
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").Value = "'7.05"
$ws.Range("E26").Value = "  -2.78%  "

# Row 27
$ws.Range("E27").Value = "  +1.45%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.67"
$ws.Range("E28").Value = "  -4.66%  "

# Row 29
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = "  -3.69%  "

# Row 31
$ws.Range("D31").Value = "'0.0487"
$ws.Range("E31").Value = "  -3.56%  "

# Row 32
$ws.Range("D32").Value = "'3.30"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").Value = "'3.18"
$ws.Range("E33").Value = "  +1.40%  "

# Row 34
$ws.Range("D34").Value = "1.424.32"
$ws.Range("E34").Value = "  -2.58%  "

# Row 35
$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = "  +1.00%  "

# Row 36
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "  +0.29%  "

# Row 37
$ws.Range("D37").Value = "'0.568"
$ws.Range("E37").Value = "  -1.68%  "

# Row 38
$ws.Range("D38").Value = "'0.889"
$ws.Range("E38").Value = "  -4.51%  "

# Row 39
$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  -3.45%  "

# Row 40
$ws.Range("E40").Value = "  -0.96%  "

# Row 41
$ws.Range("D41").Value = "'1.00"

# Row 42
$ws.Range("D42").Value = "'0.828"
$ws.Range("E42").Value = "  +4.31%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "  +2.27%  "

# Row 44
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").Value = "'2.46"
$ws.Range("E44").Value = "  -1.82%  "

# Row 46
$ws.Range("D46").Value = "'64.84"
$ws.Range("E46").Value = "  -7.45%  "

# Row 47
$ws.Range("D47").Value = "1.787.20"

# Row 48
$ws.Range("D48").Value = "'1.69"
$ws.Range("E48").Value = "  -3.68%  "

# Row 49
$ws.Range("D49").Value = "'88.42"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  -0.27%  "

# Row 51
$ws.Range("D51").Value = "'0.0995"
$ws.Range("E51").Value = "  -3.13%  "
